# Add translations for steel industry plot
#
# 1. Remove the stale translator note comment on B100.
# 2. Tweak one existing translation (red-box -> box wording).
# 3. Append 14 new EN/ZH label pairs for the steel-industry weekly
#    operating indicators chart, as rows 102-115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old reviewer comment that lived on B100.
$ws.Range("B100").Comment.Delete()

# 2) Fix existing translation wording.
$ws.Range("B101").Value = "方框内显示当年同比变化"

# 3) New translation rows.
$pairs = @(
    @("Steel industry weekly operating indicators", "钢铁行业周运行数据"),
    @("Average Daily Output: Crude Steel", "平均日产量：粗钢"),
    @("Blast furnace capacity utilization", "高炉产能利用率"),
    @("Blast furnace starting rate", "高炉开工率"),
    @("Capacity utilization: Electric Stove", "产能利用率：电炉"),
    @("Deformed Steel Bar: Operating Rate", "螺纹钢：开工率"),
    @("Estimated Average Daily Output: Crude Steel", "平均估算日产量：粗钢"),
    @("Estimated Daily Average Output: Pig Iron", "平均估算日产量：生铁"),
    @("Estimated Daily Average Output: Steel Products", "平均估算日产量：钢铁产品"),
    @("Operating Rate of Blast Furnaces: Tangshan", "高炉开工率：唐山"),
    @("Operating Rate: Electric Furnace", "开工率：电炉"),
    @("Tangshan: Operating Rate of Blast Furnaces", "唐山：高炉开工率"),
    @("Wire Rod: Operating Rate of Main Steel Plant", "线材：主要钢厂开工率"),
    @("Source: Wind Information", "数据来源：万得资讯")
)

$row = 102
foreach ($pair in $pairs) {
    $ws.Range("A$row").Value = $pair[0]
    $ws.Range("B$row").Value = $pair[1]
    $row = $row + 1
}

# Match the existing label-table formatting (style index reused, not a
# fresh one) by copying the format of the last pre-existing data row down
# across the whole freshly written block in one shot.
$ws.Range("A101:B101").Copy()
$ws.Range("A102:B115").PasteSpecial(-4122)
